# This script reproduces a stock re-sort / quantity-correction pass over the
# "CryCompanywiseStockReport" sheet: a number of duplicate-SKU rows were
# re-ordered (their Code/Rate/Qty/Value columns were moved to the correct
# row), a handful of Qty (F) values were reduced by a few units with Value (G)
# recalculated as Rate * Qty, and every Sub Total / Grand Total cell (column B)
# was refreshed to reflect the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 106, 107 ---
$ws.Cells.Item(106,6).Value = 129  # F106
$ws.Cells.Item(106,7).Value = 19436.43  # G106
$ws.Cells.Item(107,6).Value = 68  # F107
$ws.Cells.Item(107,7).Value = 6967.28  # G107

# --- Rows 114 ---
$ws.Cells.Item(114,2).Value = 263450.86  # B114

# --- Rows 146, 147, 148, 150, 152 ---
$ws.Cells.Item(146,2).Value = 53925  # B146
$ws.Cells.Item(146,5).Value = 79.37  # E146
$ws.Cells.Item(146,6).Value = 1  # F146
$ws.Cells.Item(146,7).Value = 66.44  # G146
$ws.Cells.Item(147,2).Value = 64350  # B147
$ws.Cells.Item(147,5).Value = 70.63  # E147
$ws.Cells.Item(147,6).Value = 2  # F147
$ws.Cells.Item(147,7).Value = 132.88  # G147
$ws.Cells.Item(148,2).Value = 57756  # B148
$ws.Cells.Item(148,6).Value = -100  # F148
$ws.Cells.Item(148,7).Value = -6644  # G148
$ws.Cells.Item(150,6).Value = 279  # F150
$ws.Cells.Item(150,7).Value = 27009.99  # G150
$ws.Cells.Item(152,2).Value = 21974.75  # B152

# --- Rows 197, 200 ---
$ws.Cells.Item(197,6).Value = 73  # F197
$ws.Cells.Item(197,7).Value = 3393.77  # G197
$ws.Cells.Item(200,2).Value = 47877.84  # B200

# --- Rows 233, 234 ---
$ws.Cells.Item(233,2).Value = 48719  # B233
$ws.Cells.Item(233,5).Value = 353.35  # E233
$ws.Cells.Item(233,6).Value = -81  # F233
$ws.Cells.Item(233,7).Value = -23955.75  # G233
$ws.Cells.Item(234,2).Value = 64979  # B234
$ws.Cells.Item(234,5).Value = 314.41  # E234
$ws.Cells.Item(234,6).Value = 5  # F234
$ws.Cells.Item(234,7).Value = 1478.75  # G234

# --- Rows 255 ---
$ws.Cells.Item(255,6).Value = 52  # F255
$ws.Cells.Item(255,7).Value = 1667.64  # G255

# --- Rows 274 ---
$ws.Cells.Item(274,2).Value = 86465.42999999999  # B274

# --- Rows 292, 293, 294, 295 ---
$ws.Cells.Item(292,2).Value = 63520  # B292
$ws.Cells.Item(292,5).Value = 153.4  # E292
$ws.Cells.Item(292,6).Value = 73  # F292
$ws.Cells.Item(292,7).Value = 10532.44  # G292
$ws.Cells.Item(293,2).Value = 55373  # B293
$ws.Cells.Item(293,5).Value = 163.62  # E293
$ws.Cells.Item(293,6).Value = -94  # F293
$ws.Cells.Item(293,7).Value = -13562.32  # G293
$ws.Cells.Item(294,2).Value = 57802  # B294
$ws.Cells.Item(294,5).Value = 162.71  # E294
$ws.Cells.Item(294,6).Value = -79  # F294
$ws.Cells.Item(294,7).Value = -11334.92  # G294
$ws.Cells.Item(295,2).Value = 63571  # B295
$ws.Cells.Item(295,5).Value = 152.53  # E295
$ws.Cells.Item(295,6).Value = 4  # F295
$ws.Cells.Item(295,7).Value = 573.92  # G295

# --- Rows 299, 300 ---
$ws.Cells.Item(299,2).Value = 63510  # B299
$ws.Cells.Item(299,5).Value = 50.66  # E299
$ws.Cells.Item(299,6).Value = 145  # F299
$ws.Cells.Item(299,7).Value = 6907.8  # G299
$ws.Cells.Item(300,2).Value = 55356  # B300
$ws.Cells.Item(300,5).Value = 54.04  # E300
$ws.Cells.Item(300,6).Value = -158  # F300
$ws.Cells.Item(300,7).Value = -7527.12  # G300

# --- Rows 311, 312 ---
$ws.Cells.Item(311,2).Value = 61605  # B311
$ws.Cells.Item(311,5).Value = 133.78  # E311
$ws.Cells.Item(311,6).Value = -13  # F311
$ws.Cells.Item(311,7).Value = -1455.48  # G311
$ws.Cells.Item(312,2).Value = 63563  # B312
$ws.Cells.Item(312,5).Value = 119.04  # E312
$ws.Cells.Item(312,6).Value = 2  # F312
$ws.Cells.Item(312,7).Value = 223.92  # G312

# --- Rows 328 ---
$ws.Cells.Item(328,6).Value = 783  # F328
$ws.Cells.Item(328,7).Value = 16466.49  # G328

# --- Rows 333 ---
$ws.Cells.Item(333,6).Value = 714  # F333
$ws.Cells.Item(333,7).Value = 122329.62  # G333

# --- Rows 339 ---
$ws.Cells.Item(339,2).Value = 311087.17  # B339

# --- Rows 356, 357 ---
$ws.Cells.Item(356,2).Value = 63681  # B356
$ws.Cells.Item(356,5).Value = 23.84  # E356
$ws.Cells.Item(356,6).Value = 0  # F356
$ws.Cells.Item(356,7).Value = 0  # G356
$ws.Cells.Item(357,2).Value = 31930  # B357
$ws.Cells.Item(357,5).Value = 26.8  # E357
$ws.Cells.Item(357,6).Value = -62  # F357
$ws.Cells.Item(357,7).Value = -1390.04  # G357

# --- Rows 420, 421 ---
$ws.Cells.Item(420,2).Value = 47097  # B420
$ws.Cells.Item(420,4).Value = 112.28  # D420
$ws.Cells.Item(420,5).Value = 134.16  # E420
$ws.Cells.Item(420,6).Value = 15  # F420
$ws.Cells.Item(420,7).Value = 1684.2  # G420
$ws.Cells.Item(421,2).Value = 58047  # B421
$ws.Cells.Item(421,4).Value = 105.54  # D421
$ws.Cells.Item(421,5).Value = 126.1  # E421
$ws.Cells.Item(421,6).Value = 42  # F421
$ws.Cells.Item(421,7).Value = 4432.68  # G421

# --- Rows 446, 448 ---
$ws.Cells.Item(446,6).Value = 144  # F446
$ws.Cells.Item(446,7).Value = 8843.040000000001  # G446
$ws.Cells.Item(448,2).Value = 39900  # B448

# --- Rows 465, 466, 467, 468 ---
$ws.Cells.Item(465,2).Value = 65069  # B465
$ws.Cells.Item(465,5).Value = 14.3  # E465
$ws.Cells.Item(465,6).Value = 2  # F465
$ws.Cells.Item(465,7).Value = 26.9  # G465
$ws.Cells.Item(466,2).Value = 53757  # B466
$ws.Cells.Item(466,5).Value = 16.08  # E466
$ws.Cells.Item(466,6).Value = -159  # F466
$ws.Cells.Item(466,7).Value = -2138.55  # G466
$ws.Cells.Item(467,2).Value = 65068  # B467
$ws.Cells.Item(467,5).Value = 13.97  # E467
$ws.Cells.Item(467,6).Value = 111  # F467
$ws.Cells.Item(467,7).Value = 1459.65  # G467
$ws.Cells.Item(468,2).Value = 53602  # B468
$ws.Cells.Item(468,5).Value = 15.69  # E468
$ws.Cells.Item(468,6).Value = -231  # F468
$ws.Cells.Item(468,7).Value = -3037.65  # G468

# --- Rows 472, 473 ---
$ws.Cells.Item(472,2).Value = 64915  # B472
$ws.Cells.Item(472,5).Value = 20.98  # E472
$ws.Cells.Item(472,6).Value = 0  # F472
$ws.Cells.Item(472,7).Value = 0  # G472
$ws.Cells.Item(473,2).Value = 45695  # B473
$ws.Cells.Item(473,5).Value = 23.58  # E473
$ws.Cells.Item(473,6).Value = -36  # F473
$ws.Cells.Item(473,7).Value = -710.28  # G473

# --- Rows 479, 480 ---
$ws.Cells.Item(479,2).Value = 45718  # B479
$ws.Cells.Item(479,5).Value = 19.38  # E479
$ws.Cells.Item(479,6).Value = -294  # F479
$ws.Cells.Item(479,7).Value = -4768.68  # G479
$ws.Cells.Item(480,2).Value = 64927  # B480
$ws.Cells.Item(480,5).Value = 17.26  # E480
$ws.Cells.Item(480,6).Value = 183  # F480
$ws.Cells.Item(480,7).Value = 2968.26  # G480

# --- Rows 490, 491, 492 ---
$ws.Cells.Item(490,2).Value = 65067  # B490
$ws.Cells.Item(490,5).Value = 15.65  # E490
$ws.Cells.Item(490,6).Value = 246  # F490
$ws.Cells.Item(490,7).Value = 3623.58  # G490
$ws.Cells.Item(491,2).Value = 53595  # B491
$ws.Cells.Item(491,5).Value = 17.61  # E491
$ws.Cells.Item(491,6).Value = -335  # F491
$ws.Cells.Item(491,7).Value = -4934.55  # G491
$ws.Cells.Item(492,2).Value = -3548.72  # B492

# --- Rows 559 ---
$ws.Cells.Item(559,6).Value = 103  # F559
$ws.Cells.Item(559,7).Value = 2044.55  # G559

# --- Rows 564 ---
$ws.Cells.Item(564,2).Value = 8101.34  # B564

# --- Rows 576, 577 ---
$ws.Cells.Item(576,2).Value = 64810  # B576
$ws.Cells.Item(576,5).Value = 291.22  # E576
$ws.Cells.Item(576,6).Value = 6  # F576
$ws.Cells.Item(576,7).Value = 1643.52  # G576
$ws.Cells.Item(577,2).Value = 53319  # B577
$ws.Cells.Item(577,5).Value = 310.64  # E577
$ws.Cells.Item(577,6).Value = -6  # F577
$ws.Cells.Item(577,7).Value = -1643.52  # G577

# --- Rows 608, 609 ---
$ws.Cells.Item(608,2).Value = 64830  # B608
$ws.Cells.Item(608,5).Value = 34.9  # E608
$ws.Cells.Item(608,6).Value = 112  # F608
$ws.Cells.Item(608,7).Value = 3676.96  # G608
$ws.Cells.Item(609,2).Value = 60022  # B609
$ws.Cells.Item(609,5).Value = 37.22  # E609
$ws.Cells.Item(609,6).Value = -113  # F609
$ws.Cells.Item(609,7).Value = -3709.79  # G609

# --- Rows 714, 717, 718 ---
$ws.Cells.Item(714,6).Value = 38  # F714
$ws.Cells.Item(714,7).Value = 3099.28  # G714
$ws.Cells.Item(717,2).Value = 63150  # B717
$ws.Cells.Item(717,4).Value = 75.68000000000001  # D717
$ws.Cells.Item(717,5).Value = 80.45  # E717
$ws.Cells.Item(717,6).Value = 64  # F717
$ws.Cells.Item(717,7).Value = 4843.52  # G717
$ws.Cells.Item(718,2).Value = 61428  # B718
$ws.Cells.Item(718,4).Value = 69.16  # D718
$ws.Cells.Item(718,5).Value = 73.52  # E718
$ws.Cells.Item(718,6).Value = 1  # F718
$ws.Cells.Item(718,7).Value = 69.16  # G718

# --- Rows 728 ---
$ws.Cells.Item(728,2).Value = 140781.79  # B728

# --- Rows 801, 804, 805, 806 ---
$ws.Cells.Item(801,6).Value = 35  # F801
$ws.Cells.Item(801,7).Value = 1362.55  # G801
$ws.Cells.Item(804,2).Value = 78829.12  # B804
$ws.Cells.Item(805,2).Value = 3062112.1  # B805
$ws.Cells.Item(806,2).Value = 3062112.1  # B806
